# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets, which contain identical data tables.

$wb = $excel.ActiveWorkbook

$changes = @{
    2  = 2849
    4  = 97
    5  = 6690
    6  = 1585
    7  = 16
    9  = 50
    10 = 103
    11 = 18
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Range("F$row").Value = $changes[$row]
    }
}
